$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new default exercise ("Squat") description, sets, and reps on row 24,
# next to the existing "Squat" name cell in column A.
$ws.Range("B24").Value = "Begin by positioning the barbell on the squat rack at a height where you can comfortably lift it on the pads of your shoulders while keeping your feet planted on the floor. Stand with your feet at shoulder width. Lift the weight off with the bar on your upper back, holding it with your hands in a position that feels easy to keep stable. Squat down by moving your hips back, trying to keep your core tense and back from rounding. You should squat until your upper leg is approximately parallel with the floor. Return to a standing position with a similar motion and repeat."
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 5

# Match the reported selection/view from the edit session.
$ws.Range("E29").Select()
